$d = $word.ActiveDocument
$dash = [char]0x2013

# ------------------------------------------------------------------
# Edit 1: the "refund" paragraph was split across 4 runs purely because
# of inline formatting marks around "DOES NOT" - collapse it back down
# to a single run carrying the full sentence.
# ------------------------------------------------------------------
$startRange = $d.Content
$startRange.Find.Execute("it is assumed that for a refund to be called", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$startPos = $startRange.Start

$endAnchor = $d.Content
$endAnchor.Find.Execute("as both functions are indepenent tasks", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$lastRunStart = $endAnchor.Start

$mergedLead = "it is assumed that for a refund to be called $dash a capture DOES NOT called needs to be carried beforehand "

$middleRange = $d.Range($startPos, $lastRunStart)
$middleRange.Delete()

$insertPoint = $d.Range($startPos, $startPos)
$insertPoint.InsertBefore($mergedLead)

# ------------------------------------------------------------------
# Edit 2: fill in the trailing empty paragraph (just before the
# section break) with the new assumption about the transaction DB.
# ------------------------------------------------------------------
$lastParagraph = $d.Paragraphs.Item($d.Paragraphs.Count)
$lastRange = $lastParagraph.Range
$lastRange.MoveEnd(1, -1) | Out-Null
$lastRange.InsertAfter("since we do not have a valid transaction ID at the offline verification step $dash these checks cannot be included in the transaction database")

# ------------------------------------------------------------------
# Edit 3: register the new "ListLabel 5" character style (mirrors the
# existing ListLabel1..4 styles) used by the renumbered list markers.
# ------------------------------------------------------------------
$newStyle = $d.Styles.Add("ListLabel5", 2)
$newStyle.NameLocal = "ListLabel 5"
$newStyle.QuickStyle = $true
$newStyle.Font.NameBi = "OpenSymbol"

Write-Output "done"
